$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ВСОШ по физре"
$ws.Range("B2").Value = "Безруков Владислав Александрович"
$ws.Range("D2").Value = "Технология"
$ws.Range("E2").Value = "9в"
$ws.Range("F2").Value = "Спорт"
$ws.Range("G2").Value = "Школьный"
$ws.Range("H2").Value = "Призёр"
$ws.Range("J2").Value = "22.07.2021"
$ws.Range("K2").Value = "26.02.2023"
$ws.Range("L2").Value = "26.02.2023"

# Row 3
$ws.Range("A3").Value = "ВСОШ по информатике"
$ws.Range("B3").Value = "Жинжило Татьяна Кирилловна"
$ws.Range("D3").Value = "Русский язык"
$ws.Range("E3").Value = "9в"
$ws.Range("F3").Value = "Спорт"
$ws.Range("G3").Value = "Школьный"

# "12.12.2012" would otherwise be auto-recognised as a date by Excel's
# smart entry, so force text formatting first, then strip the formatting
# back off so no stray style is left behind on the cell.
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "12.12.2012"
$ws.Range("J3").ClearFormats()

$ws.Range("K3").Value = "26.02.2023"
$ws.Range("L3").Value = "26.02.2023"

# Remove the "ФИО преподавателя" values in column C for rows 2 and 3 only
# (header C1 stays, and D..L stay in place - not shifted)
$ws.Range("C2:C3").ClearContents()
